$d = $word.ActiveDocument

$replacements = @(
    @("96÷7=13, 5", "45÷9=5, 0"),
    @("99÷7=14, 1", "69÷8=8, 5"),
    @("82÷2=41, 0", "95÷6=15, 5"),
    @("95÷4=23, 3", "31÷8=3, 7"),
    @("24÷9=2, 6", "11÷8=1, 3"),
    @("77÷9=8, 5", "34÷5=6, 4"),
    @("86÷9=9, 5", "34÷8=4, 2"),
    @("42÷2=21, 0", "74÷3=24, 2"),
    @("90÷7=12, 6", "70÷4=17, 2"),
    @("14÷8=1, 6", "71÷9=7, 8"),
    @("18÷6=3, 0", "57÷6=9, 3"),
    @("88÷2=44, 0", "56÷2=28, 0"),
    @("58÷2=29, 0", "45÷9=5, 0"),
    @("37÷8=4, 5", "41÷7=5, 6"),
    @("78÷6=13, 0", "80÷3=26, 2"),
    @("72÷9=8, 0", "80÷4=20, 0"),
    @("92÷4=23, 0", "86÷8=10, 6"),
    @("44÷9=4, 8", "12÷8=1, 4"),
    @("73÷2=36, 1", "38÷4=9, 2"),
    @("95÷9=10, 5", "26÷8=3, 2"),
    @("24÷5=4, 4", "94÷5=18, 4"),
    @("14÷6=2, 2", "27÷6=4, 3"),
    @("56÷4=14, 0", "89÷5=17, 4"),
    @("75÷4=18, 3", "91÷6=15, 1"),
    @("70÷5=14, 0", "51÷7=7, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
